# "arreglín en el excel"
#
# The balance-sheet rows had gotten scattered diagonally across the sheet
# (row 4 shifted 1 column right, row 5 shifted 2 columns right, row 6
# shifted 5 columns right, and so on). This fixes it by sliding every
# row's cells back so the table starts at column A and only spans A:C.
#
# Single-cell Range.Cut(destination) is used instead of a Value/Value2
# assignment: it moves the cell (clearing the source) while preserving the
# original cell type/formatting exactly, even for text that merely looks
# numeric (e.g. "232.245"), without touching the workbook's shared style
# table. Each row is processed left-to-right so a source cell is always
# read before a later move in the same row could overwrite it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Move-Cell($srcAddr, $dstAddr) {
    if ($srcAddr -eq $dstAddr) { return }
    $ws.Range($srcAddr).Cut($ws.Range($dstAddr))
}

# Row 4: B4 -> A4
Move-Cell "B4" "A4"

# Row 5: C5:E5 -> A5:C5
Move-Cell "C5" "A5"
Move-Cell "D5" "B5"
Move-Cell "E5" "C5"

# Row 6: F6:H6 -> A6:C6
Move-Cell "F6" "A6"
Move-Cell "G6" "B6"
Move-Cell "H6" "C6"

# Row 7: I7:K7 -> A7:C7
Move-Cell "I7" "A7"
Move-Cell "J7" "B7"
Move-Cell "K7" "C7"

# Row 8: L8:N8 -> A8:C8
Move-Cell "L8" "A8"
Move-Cell "M8" "B8"
Move-Cell "N8" "C8"

# Row 9: O9:Q9 -> A9:C9
Move-Cell "O9" "A9"
Move-Cell "P9" "B9"
Move-Cell "Q9" "C9"

# Row 10: R10:T10 -> A10:C10
Move-Cell "R10" "A10"
Move-Cell "S10" "B10"
Move-Cell "T10" "C10"
